$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate()

$ws.Range("C1").Value = "IF more than one agent make sure that candidate power plants also have the producer"

$ws.Range("C2").Select()
